# The template has a footer paragraph holding a "legacy" Word field:
#
#     <w:r><w:fldChar w:fldCharType="begin"/></w:r>
#     <w:r><w:instrText xml:space="preserve"> m:'doc.html'.fromHTMLURI() </w:instrText></w:r>
#     <w:r><w:fldChar w:fldCharType="end"/></w:r>
#
# The parser was switched to TokenIteratorFieldRewriterSplit, which reads
# M2Doc queries straight out of plain run text instead of Word field
# codes, so the field must become a single run of literal text using the
# M2Doc bracket syntax:
#
#     <w:r><w:t xml:space="preserve">{m:'doc.html'.fromHTMLURI()}</w:t></w:r>
#
# Find the field (searching every section/footer defensively, although
# this template only has one) by inspecting each field's code for the
# fromHTMLURI() call used by this particular footer.

$d = $word.ActiveDocument

$targetFtr = $null
$targetField = $null
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)
    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $flds = $ftr.Range.Fields
            for ($k = 1; $k -le $flds.Count; $k++) {
                $candidate = $flds.Item($k)
                $code = $candidate.Code.Text
                if ($code -ne $null -and $code.Contains("fromHTMLURI")) {
                    $targetFtr = $ftr
                    $targetField = $candidate
                }
            }
        }
    }
}

if ($targetField -ne $null) {
    $newText = "{m:'doc.html'.fromHTMLURI()}"

    # Locate the paragraph that follows the field's own paragraph (by its
    # visible text) before touching anything, so we can re-find it once
    # the field is gone.
    $after = $targetFtr.Range.Duplicate
    $after.Find.Execute("End of", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

    # Deleting the field removes the begin/instrText/end runs but keeps
    # its (now empty) paragraph in place.
    $targetField.Delete()

    # Re-find the same following paragraph after the mutation and step
    # one character to the left: that lands exactly inside the emptied
    # field paragraph, which is where the replacement text belongs.
    $after2 = $targetFtr.Range.Duplicate
    $after2.Find.Execute("End of", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $targetPos = $after2.Start - 1

    $insertion = $targetFtr.Range.Duplicate
    $insertion.Start = $targetPos
    $insertion.End = $targetPos
    $insertion.InsertBefore($newText)
}
